$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Update the "Status" text from "Ready for handoff" to "In Translation"
# everywhere it appears: the Overview sheet's per-language status columns
# (zh-cn, de-de) and each language sheet's own "Status" column.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: compare with the literal string on the left-hand side so
        # boolean-valued cells (e.g. "True"/"False") aren't coerced into
        # matching a non-empty string on the right-hand side.
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
        }
    }
}

# Re-fit the columns whose text just got shorter so the stored column
# widths reflect the new content, matching Excel's behavior of
# auto-sizing a column after its text changes.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null
